$d = $word.ActiveDocument

# Step 1: change text "1480x720" -> "1280x720" universally using Find/Replace
$null = $d.Content.Find.Execute("1480x720", $false, $false, $false, $false, $false, $true, 1, $false, "1280x720", 2)

# Step 2: Split "Be sure resolution is set to 1280x720" after "12" (before "80x720")
$r = $d.Content
$null = $r.Find.Execute("Be sure resolution is set to 12", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $r.End
$sub = $d.Range($splitPoint, $splitPoint + 6)
$sub.Font.Bold = 1
$sub.Font.Bold = 0

# Step 3: Split ": settings>advanced settings>1280x720" after "12" (before "80x720")
$r2 = $d.Content
$null = $r2.Find.Execute("settings>advanced settings>12", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "match2: [$($r2.Text)] start=$($r2.Start) end=$($r2.End)"
$splitPoint2 = $r2.End
$sub2 = $d.Range($splitPoint2, $splitPoint2 + 6)
Write-Host "sub2 to toggle: [$($sub2.Text)]"
$sub2.Font.Bold = 1
$sub2.Font.Bold = 0

Write-Host "Doc text now:"
Write-Host $d.Paragraphs(4).Range.Text
